$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.65%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.34%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.589"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.20%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08191"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.02%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.753"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.50%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.374"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.91%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.883"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.48%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.805"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.17%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9429"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.71%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1189"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.67%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1914"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.38%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09695"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.33%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04351"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "11.69%"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.85%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001282"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.17%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005915"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.63%"

$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.004345"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.15%"

$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.531"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.74%"

$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3536"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.07%"

$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.730"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.74%"

$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1370"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.27%"

$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2495"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.45%"

$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04386"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.11%"

$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001240"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.50%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001236"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.84%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004007"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "31.59%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02769"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.47%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05701"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.58%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007966"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.30%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009742"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.43%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.29%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002106"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.51%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01004"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.41%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007319"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.31%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.33%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003448"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002280"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002110"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.33%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.33%"
